$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Replace the end of the "RunMe" bullet sentence:
#    ". A graphical user interface (GUI) is implemented and the user may
#    choose whether or not to use the GUI or the interactions pane of the
#    java program to run the column calculator."
#    ->
#    ". The user can then view the column calculator through the
#    interactions pane of the desired java run program."
# ---------------------------------------------------------------------
$old1 = "A graphical user interface (GUI) is implemented and the user may choose whether or not to use the GUI or the interactions pane of the java program to run the column calculator."
$new1 = "The user can then view the column calculator through the interactions pane of the desired java run program."
$r1 = $d.Content
$found1 = $r1.Find.Execute($old1, $true, $false, $false, $false, $false, $true, 1, $false, $new1, 2)

# ---------------------------------------------------------------------
# 2) Merge the three runs that spell out " equilibrium.txt and inputs.txt "
#    into a single run with identical text (no visible text change, just
#    a run-boundary cleanup). Re-typing the phrase over itself collapses
#    the three original runs into one.
# ---------------------------------------------------------------------
$old2 = " equilibrium.txt and inputs.txt "
$r2 = $d.Content
$found2 = $r2.Find.Execute($old2, $true, $false, $false, $false, $false, $true, 1, $false, $old2, 2)

# ---------------------------------------------------------------------
# 3) Merge "Thank you again for using the Group 3 column calcul" + "ator"
#    (which were split by the old _GoBack bookmark) back into one run.
# ---------------------------------------------------------------------
$old3 = "Thank you again for using the Group 3 column calculator"
$r3 = $d.Content
$found3 = $r3.Find.Execute($old3, $true, $false, $false, $false, $false, $true, 1, $false, $old3, 2)

# ---------------------------------------------------------------------
# 4) Move the "_GoBack" bookmark (Word's "last edit" marker) so it now
#    splits the word "computer" into "comput" | "er", matching the most
#    recent edit location. Adding a bookmark with a name that already
#    exists relocates it (removing it from its old position).
# ---------------------------------------------------------------------
$r4 = $d.Content
$found4 = $r4.Find.Execute("comput", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found4) {
    $d.Bookmarks.Add("_GoBack", $r4)
}

Write-Host "found1=$found1 found2=$found2 found3=$found3 found4=$found4"
